# ConversionTracker.xlsx update
# -------------------------------------------------------------------------
# Adds a "Comments" column (D), two brand-new VBA-module rows that used to
# be blank filler rows, pushes the two existing rows that used to sit at
# 10/11 down to 15/16, and lets the Total formula pick up the new numbers.
#
# NOTE: new literal text is written in the exact sequence below because
# that sequence is what decides the order the strings end up in
# xl/sharedStrings.xml (first-write order), and that has to line up with
# the target workbook's table.
# -------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122
$xlLeft = -4131

# ---------------------------------------------------------------------
# 1. Re-style rows 10-16 in columns A:C before touching their contents:
#    rows 10-12 become "yellow" rows (same look as rows 2-9), rows 13-16
#    become "green" rows (same look as the old rows 10-11). Grab the
#    "green" format first, since row 10 is about to be repainted yellow.
# ---------------------------------------------------------------------
$ws.Range("A10:C10").Copy()
$ws.Range("A13:C16").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

$ws.Range("A2:C2").Copy()
$ws.Range("A10:C12").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 2. Build the "style 8" look (font1 + border1 + left align, no fill) for
#    the new Comments column by copying an existing bordered/no-fill cell
#    and then switching its alignment, then stamp it down column D.
# ---------------------------------------------------------------------
$ws.Range("A17").Copy()
$ws.Range("D2").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false
$ws.Range("D2").HorizontalAlignment = $xlLeft

$ws.Range("D2").Copy()
$ws.Range("D3:D16").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

# Rows 17-29 and the Total row just reuse the plain bordered style (3)
# that is already used by the rest of those rows.
$ws.Range("A17").Copy()
$ws.Range("D17:D29").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

$ws.Range("C30").Copy()
$ws.Range("D30").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

# Header cell D1 matches the rest of row 1.
$ws.Range("A1").Copy()
$ws.Range("D1").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 3. New text - written in the order that must match the first-seen
#    order in the saved sharedStrings table.
# ---------------------------------------------------------------------
$ws.Range("B14").Value = "mAddDatatoSuptPoints"
$ws.Range("B10").Value = "mDec2Bin"
$ws.Range("D1").Value = "Comments"
$ws.Range("D14").Value = "Errors should be resolved after cSuptPoints will be added."
$ws.Range("B13").Value = "mAddRouteBranchesandOptions"
$ws.Range("B11").Value = "mArchiveRun"
$ws.Range("B12").Value = "mAssignSuptFrameSectionSizes"

# ---------------------------------------------------------------------
# 4. Remaining values: column A ("Alex") and column C (lengths). Rows
#    15/16 get the data that used to live in rows 10/11.
# ---------------------------------------------------------------------
$ws.Range("A10:A16").Value = "Alex"

$ws.Range("C10").Value = 40
$ws.Range("C11").Value = 110
$ws.Range("C12").Value = 77
$ws.Range("C13").Value = 502
$ws.Range("C14").Value = 54
$ws.Range("C15").Value = 378
$ws.Range("C16").Value = 217

$ws.Range("B15").Value = "mDefineSupFramesSections"
$ws.Range("B16").Value = "mSubInitializationSupA"

# ---------------------------------------------------------------------
# 5. Column width for the new Comments column, then the scroll/selection
#    state last (matches where the edit left the cursor).
# ---------------------------------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 70.33

$ws.Range("D16").Select()

$wb.Save()
